$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bold/fill/border styling that used to be on the header cells,
# restoring them to the plain "Normal" style.
$ws.Range("A1:B1").Style = "Normal"

# Add the new header columns. Values are written in this particular order
# (D, E, F, then C) so that the shared-string table ends up ordered
# WR Time, EX Time, Total Time, Needs Imaging - matching how the fields
# were actually added to the report.
$ws.Range("D1").Value = "WR Time"
$ws.Range("E1").Value = "EX Time"
$ws.Range("F1").Value = "Total Time"
$ws.Range("C1").Value = "Needs Imaging"

# "Needs Imaging" is the widest header, so give its column an autofit width.
$ws.Columns("C:C").AutoFit()

# Leave the selection where the author left off editing.
$ws.Range("C5").Select()
